$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single "7118 Almaden Place, San Jose CA 95120" address line
#    (the one in the letter header, not the one inside the table further down)
#    into two separate paragraphs: "7118 Almaden Place" and "San Jose, CA 95120".
$addressParaIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "7118 Almaden Place, San Jose CA 95120*" -and $p.Range.Information(12) -eq $false) {
        $addressParaIndex = $i
    }
}

if ($addressParaIndex -ne -1) {
    $addrPara = $d.Paragraphs.Item($addressParaIndex)
    $addrRange = $addrPara.Range
    $addrRange.Find.Execute(", San Jose CA 95120", $true, $false, $false, $false, $false, $true, 1, $false, "^pSan Jose, CA 95120", 2) | Out-Null

    # Re-apply the run-level character formatting (Arial 11pt, incl. complex-script)
    # to the newly created "San Jose, CA 95120" paragraph, matching its sibling runs.
    $newParaIndex = $addressParaIndex + 1
    $newPara = $d.Paragraphs.Item($newParaIndex)
    $newFont = $newPara.Range.Font
    $newFont.NameAscii = "Arial"
    $newFont.NameOther = "Arial"
    $newFont.NameBi = "Arial"
    $newFont.Size = 11
    $newFont.SizeBi = 11
}

# 3. Remove the now-redundant empty "No Spacing" paragraph that sat directly
#    beneath "Board of Directors".
$boardParaIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Board of Directors*") {
        $boardParaIndex = $i
    }
}

if ($boardParaIndex -ne -1) {
    $emptyParaIndex = $boardParaIndex + 1
    $emptyPara = $d.Paragraphs.Item($emptyParaIndex)
    if ($emptyPara.Range.Text.Trim() -eq "") {
        $emptyPara.Range.Delete()
    }
}
